$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 352 (shifts old rows 352-368 down to 356-372)
$ws.Range("A352:T355").Insert()

# Common / constant values for all 4 new rows (same as the surrounding template rows)
$A = 2
$B = 'Comercializadora del Agro de Limarí'
$C = 'Coquimbo'
$D = 44615
$E = 4
$F = 'Fruta'
$G = 100106
$H = 'Oleaginosos'
$I = 100106002
$J = 'Palta'
$R = 'Provincia de Limarí'
$T = 1

$rows = @(
    @{ Row = 352; K = 'Hass'; L = 'Especial'; M = 300; N = 2400; O = 2500; P = 2450; Q = '$/kilo (en caja de 17 kilos)'; S = 2450 },
    @{ Row = 353; K = 'Hass'; L = 'Primera';  M = 240; N = 2100; O = 2200; P = 2150; Q = '$/kilo (en caja de 17 kilos)'; S = 2150 },
    @{ Row = 354; K = 'Hass'; L = 'Segunda';  M = 240; N = 1800; O = 1900; P = 1850; Q = '$/kilo (en caja de 17 kilos)'; S = 1850 },
    @{ Row = 355; K = 'Hass'; L = 'Tercera';  M = 200; N = 1300; O = 1400; P = 1350; Q = '$/kilo (en caja de 17 kilos)'; S = 1350 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $A
    $ws.Cells.Item($row, 2).Value = $B
    $ws.Cells.Item($row, 3).Value = $C
    $ws.Cells.Item($row, 4).Value = $D
    $ws.Cells.Item($row, 5).Value = $E
    $ws.Cells.Item($row, 6).Value = $F
    $ws.Cells.Item($row, 7).Value = $G
    $ws.Cells.Item($row, 8).Value = $H
    $ws.Cells.Item($row, 9).Value = $I
    $ws.Cells.Item($row, 10).Value = $J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $T
}
